# deltahouronactivity.xlsx update
# - Move the report date forward (B1: 25/03/2023 -> 07/04/2023), kept as literal text
# - Remove the "NOT" category row (row 8), shifting the rows below it up by one
# - Eliminate the trailing initialization hour-index (the last numeric row, value 10),
#   so the hour index now runs 0..9 (aligned, starting the chart index from 9 rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header date, writing it through a scratch cell + PasteSpecial(values)
# so the date-like text is stored as literal text without Excel re-inferring it as
# a date serial number and without forking B1's existing style.
$scratch = $ws.Range("Z1")
$scratch.Value = "'07/04/2023"
$scratch.Copy()
$ws.Range("B1").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# Delete the "NOT" row entirely - everything below shifts up one row
$ws.Rows.Item(8).Delete()

# Eliminate the initialization row that's now left dangling at the bottom
# (originally row 23, now row 22 after the shift above)
$ws.Rows.Item(22).Delete()
